$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.06"
$ws.Range("E2").Value = "'1.71%"
$ws.Range("D3").Value = "'27.24"
$ws.Range("E3").Value = "'1.06%"
$ws.Range("D4").Value = "'4.726"
$ws.Range("E4").Value = "'5.30%"
$ws.Range("D5").Value = "'0.06082"
$ws.Range("E5").Value = "'3.18%"
$ws.Range("D6").Value = "'6.672"
$ws.Range("E6").Value = "'0.98%"
$ws.Range("D7").Value = "'0.8466"
$ws.Range("E7").Value = "'-0.40%"
$ws.Range("D8").Value = "'0.9218"
$ws.Range("E8").Value = "'-0.51%"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'1.87%"
$ws.Range("D10").Value = "'0.05039"
$ws.Range("E10").Value = "'17.49%"
$ws.Range("D11").Value = "'0.07102"
$ws.Range("E11").Value = "'1.41%"
$ws.Range("D12").Value = "'0.03133"
$ws.Range("E12").Value = "'2.71%"
$ws.Range("D13").Value = "'0.09072"
$ws.Range("E13").Value = "'-0.27%"
$ws.Range("D14").Value = "'0.001540"
$ws.Range("E14").Value = "'0.17%"
$ws.Range("D15").Value = "'0.0006107"
$ws.Range("E15").Value = "'0.65%"
$ws.Range("D16").Value = "'0.006147"
$ws.Range("E16").Value = "'1.02%"
$ws.Range("D17").Value = "'3.453"
$ws.Range("E17").Value = "'-0.50%"
$ws.Range("D18").Value = "'3.147"
$ws.Range("E18").Value = "'-0.75%"
$ws.Range("D19").Value = "'2.168"
$ws.Range("E19").Value = "'-1.21%"
$ws.Range("D20").Value = "'0.3128"
$ws.Range("E20").Value = "'3.23%"
$ws.Range("E21").Value = "'0.91%"
$ws.Range("D22").Value = "'4.106"
$ws.Range("E22").Value = "'4.92%"
$ws.Range("D23").Value = "'0.04242"
$ws.Range("E23").Value = "'-0.36%"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'0.28%"
$ws.Range("E25").Value = "'-9.08%"
$ws.Range("E26").Value = "'0.11%"
$ws.Range("D27").Value = "'0.0001575"
$ws.Range("E27").Value = "'3.44%"
$ws.Range("D40").Value = "'0.03873"
$ws.Range("E40").Value = "'1.69%"
$ws.Range("E41").Value = "'1.53%"
$ws.Range("D42").Value = "'0.004110"
$ws.Range("E42").Value = "'-34.61%"
$ws.Range("E43").Value = "'21.39%"
$ws.Range("D44").Value = "'0.002221"
$ws.Range("E44").Value = "'1.03%"
$ws.Range("D45").Value = "'0.00005321"
$ws.Range("E45").Value = "'-0.25%"
$ws.Range("E46").Value = "'0.12%"
$ws.Range("B47").Value = "'CoinbaseStockToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.05451"
$ws.Range("E47").Value = "'5.08%"
$ws.Range("B48").Value = "'BOLO"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.1353"
$ws.Range("E48").Value = "'-46.44%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.12%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.12%"
